$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leave behind a hidden _xlnm._FilterDatabase defined name scoped to this sheet,
# covering the full used range (this mirrors what Excel leaves behind after a
# filter is applied/cleared on the data).
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$Q`$547")
$filterName.Visible = $false

# Update I column values for the affected rows: "Yes" -> "No"
$rows = @(3,4,7,12,14,15,18,19,20,21,22,23,24,27,28,29,30,31,32,34,36,37,38,39,40,41,42,45,46,48,62,63,65,75,76,77,78,83,84,86,87,89,94,95,98,103,105,106,109,110,111,112,113,114,115,118,119,120,121,122,123,125,127,128,129,130,131,132,133,136,137,139,153,154,156,166,167,168,169,174,175,177,178,180,185,186,189,194,196,197,200,201,202,203,204,205,206,209,210,211,212,213,214,216,218,219,220,221,222,223,224,227,228,230,244,245,247,257,258,259,260,265,266,268,269,271,276,277,280,285,287,288,291,292,293,294,295,296,297,300,301,302,303,304,305,307,309,310,311,312,313,314,315,318,319,321,335,336,338,348,349,350,351,356,357,359,360,362,367,368,371,376,378,379,382,383,384,385,386,387,388,391,392,393,394,395,396,398,400,401,402,403,404,405,406,409,410,412,426,427,429,439,440,441,442,447,448,450,451,453,458,459,462,467,469,470,473,474,475,476,477,478,479,482,483,484,485,486,487,489,491,492,493,494,495,496,497,500,501)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = "No"
}

# Update the sheet view (scroll position / selection) to match the saved view state
$aw = $excel.ActiveWindow
$aw.ScrollRow = 513
$aw.ScrollColumn = 5
$ws.Range("M13").Select() | Out-Null
